# Update LR-pair stats with newly computed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.022418
$ws.Range("H2").Value = 0.06725399999999999
$ws.Range("M2").Value = 3.112844666666666
$ws.Range("N2").Value = 9.338533999999999
$ws.Range("O2").Value = 0.0962303687181678
$ws.Range("P2").Value = 0.09623036871816783
$ws.Range("Q2").Value = 0.06978375173733331
$ws.Range("R2").Value = 0.6280537656359999
$ws.Range("S2").Value = 0.0962303687181678
$ws.Range("T2").Value = 0.09623036871816783

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.022418
$ws.Range("H3").Value = 0.06725399999999999
$ws.Range("O3").Value = 0.5562258596073433
$ws.Range("P3").Value = 0.5562258596073434
$ws.Range("Q3").Value = 0.4033604756353333
$ws.Range("R3").Value = 3.630244280717999
$ws.Range("S3").Value = 0.5562258596073433
$ws.Range("T3").Value = 0.5562258596073434

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.022418
$ws.Range("H4").Value = 0.06725399999999999
$ws.Range("M4").Value = 11.24229066666667
$ws.Range("N4").Value = 33.726872
$ws.Range("O4").Value = 0.3475437716744887
$ws.Range("P4").Value = 0.3475437716744888
$ws.Range("Q4").Value = 0.2520296721653333
$ws.Range("R4").Value = 2.268267049488
$ws.Range("S4").Value = 0.3475437716744887
$ws.Range("T4").Value = 0.3475437716744888
